# "added CRUD on participants" - add the missing "image" column values
# for the two existing participant rows, and move the active selection
# to D4 (the cell right below the newly-filled image data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "img.jpg"
$ws.Range("D3").Value = "img.png"

$ws.Range("D4").Select()
